# Apply scheduled market-data refresh updates to Sheets/Louisoix_Profits workbook

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1178.25
$ws.Range("I18").Value = 632.2857
$ws.Range("K18").Value = 632.2857
$ws.Range("M18").Value = -348.2857

$ws.Range("H40").Value = 6000
$ws.Range("I40").Value = 6000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 6000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -5825
$ws.Range("N40").ClearContents()

$ws.Range("H64").Value = 5189.857
$ws.Range("I64").Value = 4779.8335
$ws.Range("K64").Value = 4779.8335
$ws.Range("M64").Value = -4531.8335

$ws.Range("H67").Value = 5189.857
$ws.Range("I67").Value = 4779.8335
$ws.Range("K67").Value = 4779.8335
$ws.Range("M67").Value = -3921.8335

$ws.Range("H74").Value = 12645.723
$ws.Range("I74").Value = 6635.6665
$ws.Range("J74").Value = 24665.834
$ws.Range("K74").Value = 6635.6665
$ws.Range("L74").Value = 24665.834
$ws.Range("M74").Value = -5699.6665
$ws.Range("N74").Value = -26537.834

$ws.Range("H77").Value = 12645.723
$ws.Range("I77").Value = 6635.6665
$ws.Range("J77").Value = 24665.834
$ws.Range("K77").Value = 33178.3325
$ws.Range("L77").Value = 123329.17
$ws.Range("M77").Value = -28498.3325
$ws.Range("N77").Value = -132689.17

$ws.Range("H113").Value = 6028.222
$ws.Range("J113").Value = 6866.6665
$ws.Range("L113").Value = 6866.6665
$ws.Range("N113").Value = -13374.6665

$ws.Range("H116").Value = 15677.083
$ws.Range("J116").Value = 17312.6
$ws.Range("L116").Value = 17312.6
$ws.Range("N116").Value = -24196.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2110.138
$ws.Range("I61").Value = 2122.8462
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 2122.8462
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1910.8462
$ws.Range("N61").Value = -2424

$ws.Range("H74").Value = 1897.1904
$ws.Range("J74").Value = 2487.5557
$ws.Range("L74").Value = 2487.5557
$ws.Range("N74").Value = -4235.5557

$ws.Range("H77").Value = 1897.1904
$ws.Range("J77").Value = 2487.5557
$ws.Range("L77").Value = 12437.7785
$ws.Range("N77").Value = -21173.7785

$ws.Range("H122").Value = 1586.6552
$ws.Range("I122").Value = 1404.3462
$ws.Range("K122").Value = 4213.0386
$ws.Range("M122").Value = -1763.0386

$ws.Range("H136").Value = 2110.138
$ws.Range("I136").Value = 2122.8462
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 6368.5386
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -3818.5386
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 57499.5
$ws.Range("J76").Value = 57499.5
$ws.Range("L76").Value = 57499.5
$ws.Range("N76").Value = -58129.5

$ws.Range("H79").Value = 57499.5
$ws.Range("J79").Value = 57499.5
$ws.Range("L79").Value = 57499.5
$ws.Range("N79").Value = -59683.5

$ws.Range("H86").Value = 4079
$ws.Range("J86").Value = 4935.3076
$ws.Range("L86").Value = 4935.3076
$ws.Range("N86").Value = -7181.3076

$ws.Range("H89").Value = 4079
$ws.Range("J89").Value = 4935.3076
$ws.Range("L89").Value = 24676.538
$ws.Range("N89").Value = -35908.538

$ws.Range("H134").Value = 2017.2903
$ws.Range("I134").Value = 2017.2903
$ws.Range("K134").Value = 6051.8709
$ws.Range("M134").Value = -3516.8709

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 5500
$ws.Range("I125").Value = 5500
$ws.Range("K125").Value = 16500
$ws.Range("M125").Value = -11580

$ws.Range("H129").Value = 6308.2915
$ws.Range("I129").Value = 9727.083000000001
$ws.Range("J129").Value = 2889.5
$ws.Range("K129").Value = 29181.249
$ws.Range("L129").Value = 8668.5
$ws.Range("M129").Value = -24181.249
$ws.Range("N129").Value = -18668.5

$ws.Range("H131").Value = 10752.4
$ws.Range("J131").Value = 8162.6924
$ws.Range("L131").Value = 24488.0772
$ws.Range("N131").Value = -34568.0772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 39999.5
$ws.Range("J63").Value = 39999.5
$ws.Range("L63").Value = 39999.5
$ws.Range("N63").Value = -41371.5

$ws.Range("H66").Value = 39999.5
$ws.Range("J66").Value = 39999.5
$ws.Range("L66").Value = 119998.5
$ws.Range("N66").Value = -126862.5

$ws.Range("H121").Value = 31700
$ws.Range("J121").Value = 31700
$ws.Range("L121").Value = 31700
$ws.Range("N121").Value = -35194

$ws.Range("H132").Value = 168953.17
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2324.2
$ws.Range("I93").Value = 1629.6
$ws.Range("J93").Value = 3018.8
$ws.Range("K93").Value = 1629.6
$ws.Range("L93").Value = 3018.8
$ws.Range("M93").Value = -381.5999999999999
$ws.Range("N93").Value = -5514.8

$ws.Range("H100").Value = 4682.4
$ws.Range("J100").Value = 4996.3335
$ws.Range("L100").Value = 4996.3335
$ws.Range("N100").Value = -6078.3335

$ws.Range("H122").Value = 4680.5
$ws.Range("I122").Value = 3499
$ws.Range("J122").Value = 4916.8
$ws.Range("K122").Value = 10497
$ws.Range("L122").Value = 14750.4
$ws.Range("M122").Value = -8047
$ws.Range("N122").Value = -19650.4

$ws.Range("H138").Value = 54750
$ws.Range("I138").Value = 35000
$ws.Range("J138").Value = 61333.332
$ws.Range("K138").Value = 35000
$ws.Range("L138").Value = 61333.332
$ws.Range("M138").Value = -29860
$ws.Range("N138").Value = -71613.33199999999

$ws.Range("H140").Value = 29498.666
$ws.Range("J140").Value = 34498
$ws.Range("L140").Value = 34498
$ws.Range("N140").Value = -44858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 48681.4
$ws.Range("J46").Value = 48681.4
$ws.Range("L46").Value = 48681.4
$ws.Range("N46").Value = -49143.4

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H134").Value = 48681.4
$ws.Range("J134").Value = 48681.4
$ws.Range("L134").Value = 146044.2
$ws.Range("N134").Value = -151114.2

Write-Host "edit complete"
